# database sanity check in GUI and core
# Rewrites the sample "Row/Reference Number/..." table on Sheet1 with a
# cleaned up, sanity-checked data set (rows now keyed "row1".."row6" with
# long zero-padded serial numbers), fixes a couple of reference numbers
# and dates, and appends two brand new rows (6 and 7). Sheet2 merely
# mirrors Sheet1's "Description" strings, so no direct edits are needed
# there - it follows automatically once the shared strings move.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 2 --------------------------------------------------------------
$ws.Cells.Item(2, 2).Value = 24234
$ws.Cells.Item(2, 3).Value = "row1"
$ws.Cells.Item(2, 4).Value = "AA0000000000000000000000000100"
$ws.Cells.Item(2, 5).Value = "AA0000000000000000000000000200"
$ws.Cells.Item(2, 6).Value = 41092

# ---- Row 3 --------------------------------------------------------------
$ws.Cells.Item(3, 2).Value = 54365
$ws.Cells.Item(3, 3).Value = "row2"
$ws.Cells.Item(3, 4).Value = "AB0000000000000000000000000500"
$ws.Cells.Item(3, 5).Value = "AB0000000000000000000000000600"
$ws.Cells.Item(3, 6).Value = 41093
$ws.Cells.Item(3, 2).Font.Name = "Arial"
$ws.Cells.Item(3, 3).Font.Name = "Arial"

# ---- Row 4 --------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = 7567
$ws.Cells.Item(4, 3).Value = "row3"
$ws.Cells.Item(4, 4).Value = "AA0000000000000000000000000110"
$ws.Cells.Item(4, 5).Value = "AA0000000000000000000000000120"
$ws.Cells.Item(4, 6).Value = 41094
$ws.Cells.Item(4, 2).Font.Name = "Arial"
$ws.Cells.Item(4, 3).Font.Name = "Arial"

# ---- Row 5 --------------------------------------------------------------
$ws.Cells.Item(5, 2).Value = 578
$ws.Cells.Item(5, 3).Value = "row4"
$ws.Cells.Item(5, 4).Value = "AA0000000000000000000000000090"
$ws.Cells.Item(5, 5).Value = "AA0000000000000000000000000100"
$ws.Cells.Item(5, 6).Value = 41095
$ws.Cells.Item(5, 2).Font.Name = "Arial"
$ws.Cells.Item(5, 3).Font.Name = "Arial"

# ---- Row 6 (new) ----------------------------------------------------------
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 78
$ws.Cells.Item(6, 3).Value = "row5"
$ws.Cells.Item(6, 4).Value = "AC0000000000000000000000000090"
$ws.Cells.Item(6, 5).Value = "AC0000000000000000000000000200"
$ws.Cells.Item(6, 6).Value = 41096
$ws.Cells.Item(6, 6).NumberFormat = "MM/DD/YY"

# ---- Row 7 (new) ----------------------------------------------------------
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 87654
$ws.Cells.Item(7, 3).Value = "row6"
$ws.Cells.Item(7, 4).Value = "AB0000000000000000000000000090"
$ws.Cells.Item(7, 5).Value = "AB0000000000000000000000001100"
$ws.Cells.Item(7, 6).Value = 41097
$ws.Cells.Item(7, 6).NumberFormat = "MM/DD/YY"

# ---- Column widths for the newly-important Start/End Serial columns -----
$ws.Columns.Item(4).ColumnWidth = 33.24
$ws.Columns.Item(5).ColumnWidth = 31.3

# ---- Selection moves to the cell the author was working on --------------
$ws.Range("D2").Select()
